# Insert a new weekly record as row 83, shifting all existing rows (83:175)
# down by one (to 84:176). Excel's Insert() on an entire row shifts the
# existing row (and everything below it) down automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(83).Insert()

$ws.Cells.Item(83, 1).Value = 7
$ws.Cells.Item(83, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(83, 3).Value = "Ñuble"
$ws.Cells.Item(83, 4).Value = 45225
$ws.Cells.Item(83, 5).Value = 16
$ws.Cells.Item(83, 6).Value = 100112031
$ws.Cells.Item(83, 7).Value = "Poroto verde"
$ws.Cells.Item(83, 8).Value = "Sin especificar"
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 60
$ws.Cells.Item(83, 11).Value = 27000
$ws.Cells.Item(83, 12).Value = 28000
$ws.Cells.Item(83, 13).Value = 27500
$ws.Cells.Item(83, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(83, 15).Value = "Perú"
$ws.Cells.Item(83, 16).Value = 1100
$ws.Cells.Item(83, 17).Value = 25
$ws.Cells.Item(83, 18).Value = "Hortaliza"
